$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Multi-Core: ..." known-issue bullet (and the blank
#    paragraph with the 1080-twip indent that immediately followed
#    it) from the "Tools known issues" list.
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Multi-Core:*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $p1 = $d.Paragraphs.Item($target)
    $p2 = $d.Paragraphs.Item($target + 1)
    $r = $d.Range($p1.Range.Start, $p2.Range.End)
    $r.Delete()
}

# ------------------------------------------------------------------
# 2. The footer's cached PAGE field result drops from "6" to "5"
#    now that the document is a paragraph shorter.
# ------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$fr = $ftr.Range.Duplicate
$fr.Find.Execute("6", $true, $false, $false, $false, $false, $true, 1, $false, "5", 2)
